$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Set the new header text in U1 (adds a new shared string "发货地址")
$ws.Range("U1").Value = "发货地址"

# Update the selected cell / active cell on the sheet view to U1
$ws.Range("U1").Select()
